$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.177.64"
$ws.Range("E2").Value = "  +1.47%  "
$ws.Range("D3").Value = "3.907.98"
$ws.Range("E3").Value = "  +2.96%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'464.73"
$ws.Range("E5").Value = "  +8.50%  "
$ws.Range("D6").Value = "'144.20"
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("E7").Value = "  +2.10%  "
$ws.Range("D8").Value = "'0.998"
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("E10").Value = "  +5.98%  "
$ws.Range("D11").Value = "'0.0000343"
$ws.Range("E11").Value = "  +6.19%  "
$ws.Range("E12").Value = "  -0.90%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "'10.38"
$ws.Range("E13").Value = "  -2.46%  "
$ws.Range("B14").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C14").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D14").Value = "4.518.11"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").Value = "'15.23"
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").Value = "3.932.76"
$ws.Range("E16").Value = "  +3.91%  "
$ws.Range("E17").Value = "  -0.43%  "
$ws.Range("D18").Value = "'19.96"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("E19").Value = "  +1.99%  "
$ws.Range("D20").Value = "67.277.97"
$ws.Range("E20").Value = "  +1.34%  "
$ws.Range("D21").Value = "'432.54"
$ws.Range("E21").Value = "  +5.03%  "
$ws.Range("E22").Value = "  -4.43%  "
$ws.Range("D23").Value = "'3.32"
$ws.Range("E23").Value = "  +1.16%  "
$ws.Range("D24").Value = "'88.92"
$ws.Range("E24").Value = "  +3.72%  "
$ws.Range("D25").Value = "'38.67"
$ws.Range("E25").Value = "  +4.08%  "
$ws.Range("D26").Value = "'3.51"
$ws.Range("E26").Value = "  +6.53%  "
$ws.Range("E27").Value = "  +5.25%  "
$ws.Range("D28").Value = "'10.10"
$ws.Range("E28").Value = "  +2.16%  "
$ws.Range("D29").Value = "'9.63"
$ws.Range("E29").Value = "  -1.19%  "
$ws.Range("D30").Value = "'739.35"
$ws.Range("E30").Value = "  +5.47%  "
$ws.Range("E31").Value = "  -2.18%  "
$ws.Range("E32").Value = "  +1.83%  "
$ws.Range("E33").Value = "  -0.17%  "
$ws.Range("D34").Value = "'43.08"
$ws.Range("E34").Value = "  +5.97%  "
$ws.Range("E35").Value = "  +3.94%  "
$ws.Range("D36").Value = "'58.15"
$ws.Range("E36").Value = "  +2.95%  "
$ws.Range("E37").Value = "  -0.10%  "
$ws.Range("D38").Value = "0.0₃0791"
$ws.Range("E38").Value = "  +15.75%  "
$ws.Range("D39").Value = "'5.37"
$ws.Range("E39").Value = "  -6.53%  "
$ws.Range("D40").Value = "'3.21"
$ws.Range("E40").Value = "  +12.26%  "
$ws.Range("E41").Value = "  -0.15%  "
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").Value = "'0.140"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D44").Value = "'0.335"
$ws.Range("E44").Value = "  +3.85%  "
$ws.Range("D45").Value = "'2.79"
$ws.Range("E45").Value = "  +5.44%  "
$ws.Range("E46").Value = "  +4.42%  "
$ws.Range("D47").Value = "'3.40"
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("D48").Value = "'2.48"
$ws.Range("E48").Value = "  -4.85%  "
$ws.Range("D49").Value = "'3.15"
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("D50").Value = "'2.90"
$ws.Range("E50").Value = "  +2.43%  "
$ws.Range("D51").Value = "'143.14"
$ws.Range("E51").Value = "  +0.18%  "
